$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 255
$ws.Range("F5").Value = 127
$ws.Range("F7").Value = 389
$ws.Range("F8").Value = 4786
$ws.Range("F9").Value = 4786
$ws.Range("F10").Value = 14
$ws.Range("F11").Value = 128
$ws.Range("F12").Value = 456
$ws.Range("F13").Value = 1095
$ws.Range("F14").Value = 628
$ws.Range("F15").Value = 4366
$ws.Range("F16").Value = 168
$ws.Range("F17").Value = 172
$ws.Range("F18").Value = 73
$ws.Range("F19").Value = 220
$ws.Range("F20").Value = 3513
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("F24").Value = 3170
$ws.Range("F25").Value = 138
$ws.Range("F26").Value = 129
$ws.Range("F28").Value = 155
$ws.Range("F31").Value = 82
$ws.Range("F32").Value = 63
$ws.Range("F35").Value = 129
$ws.Range("F36").Value = 5578
$ws.Range("F37").Value = 851
$ws.Range("F38").Value = 407
$ws.Range("F39").Value = 86
$ws.Range("F40").Value = 958
$ws.Range("F41").Value = 49
$ws.Range("F42").Value = 1122
$ws.Range("F43").Value = 45
$ws.Range("F44").Value = 498
$ws.Range("F46").Value = 2002
$ws.Range("F47").Value = 299
$ws.Range("F49").Value = 706
$ws.Range("F50").Value = 855

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 13
$ws.Range("F5").Value = 1
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 125
$ws.Range("F21").Value = 41
$ws.Range("F24").Value = 748

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 202

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 202
$ws.Range("F3").Value = 13
$ws.Range("F4").Value = 255
$ws.Range("F6").Value = 127
$ws.Range("F8").Value = 389
$ws.Range("F9").Value = 4786
$ws.Range("F10").Value = 4786
$ws.Range("F14").Value = 456
$ws.Range("F15").Value = 1095
$ws.Range("F16").Value = 628
$ws.Range("F17").Value = 4366
$ws.Range("F18").Value = 168
$ws.Range("F19").Value = 172
$ws.Range("F20").Value = 73
$ws.Range("F21").Value = 220
$ws.Range("F22").Value = 3513
$ws.Range("F23").Value = 3170
$ws.Range("F24").Value = 138
$ws.Range("F25").Value = 129
$ws.Range("F26").Value = 155
$ws.Range("F29").Value = 82
$ws.Range("F30").Value = 63
$ws.Range("F31").Value = 5
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 125
$ws.Range("F36").Value = 5578
$ws.Range("F38").Value = 851
$ws.Range("F39").Value = 407
$ws.Range("F42").Value = 86
$ws.Range("F43").Value = 958
$ws.Range("F44").Value = 49
$ws.Range("F45").Value = 1122
$ws.Range("F46").Value = 498
$ws.Range("F47").Value = 2002
$ws.Range("F48").Value = 299
$ws.Range("F49").Value = 706
$ws.Range("F50").Value = 855
